$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the assembly name in B11 from "Curtain wall: alu spandrel"
# to "Curtain wall: aluminum spandrel"
$ws.Range("B11").Value = "Curtain wall: aluminum spandrel"

# Reflect the resulting active selection on the sheet
$ws.Range("B11").Select()
